# Fruta / hortaliza, semanal
#
# The weekly refresh re-sorted the "Hortaliza, Terminal Hortofrutícola Agro
# Chillán - Espinaca" record set by date. Columns A,B,C,E,F,G,H,N,Q,R are
# constant for every data row, so only D (Fecha), I (Calidad), J (Volumen),
# K/L/M (Precio mínimo/máximo/promedio ponderado), O (Origen) and P (Precio
# $/Kg) actually move between rows 2-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  D=44799; I='Primera'; J=60;  K=7000; L=7000; M=7000; O='Provincia de Diguillín'; P=700 },
    @{ Row=3;  D=44211; I='Primera'; J=28;  K=8000; L=8500; M=8214; O='Región Metropolitana';   P=821 },
    @{ Row=4;  D=44798; I='Primera'; J=80;  K=7000; L=7000; M=7000; O='Provincia de Diguillín'; P=700 },
    @{ Row=5;  D=44804; I='Primera'; J=80;  K=7000; L=7500; M=7250; O='Provincia de Diguillín'; P=725 },
    @{ Row=6;  D=44806; I='Primera'; J=120; K=7000; L=7500; M=7250; O='Provincia de Diguillín'; P=725 },
    @{ Row=7;  D=44812; I='Primera'; J=60;  K=7000; L=8000; M=7500; O='Provincia de Diguillín'; P=750 },
    @{ Row=8;  D=44819; I='Primera'; J=100; K=7000; L=8000; M=7500; O='Provincia de Diguillín'; P=750 },
    @{ Row=9;  D=44817; I='Primera'; J=60;  K=7000; L=7000; M=7000; O='Provincia de Diguillín'; P=700 },
    @{ Row=10; D=44817; I='Segunda'; J=60;  K=8000; L=8000; M=8000; O='Provincia de Diguillín'; P=800 },
    @{ Row=11; D=44813; I='Primera'; J=120; K=7000; L=7500; M=7250; O='Provincia de Diguillín'; P=725 },
    @{ Row=12; D=44810; I='Primera'; J=60;  K=7000; L=8000; M=7500; O='Provincia de Diguillín'; P=750 },
    @{ Row=13; D=44790; I='Primera'; J=60;  K=8500; L=9000; M=8750; O='Región Metropolitana';   P=875 },
    @{ Row=14; D=44791; I='Primera'; J=100; K=8500; L=9000; M=8750; O='Región Metropolitana';   P=875 },
    @{ Row=15; D=44784; I='Primera'; J=100; K=8000; L=9000; M=8500; O='Región Metropolitana';   P=850 },
    @{ Row=16; D=44782; I='Primera'; J=120; K=8000; L=9000; M=8500; O='Región Metropolitana';   P=850 },
    @{ Row=17; D=44203; I='Primera'; J=27;  K=7000; L=8000; M=7556; O='Región Metropolitana';   P=756 },
    @{ Row=18; D=44775; I='Primera'; J=60;  K=8000; L=8000; M=8000; O='Región Metropolitana';   P=800 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($row, 9).Value  = $r.I   # I: Calidad
    $ws.Cells.Item($row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio mínimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio máximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
}
